$d = $word.ActiveDocument

$pairs = @(
    ,@("2025-08-20 Wednesday", "2025-08-21 Thursday")
    ,@("90-54=36", "58-50=8")
    ,@("67-66=1", "61+13=74")
    ,@("80-8=72", "87-76=11")
    ,@("85-73=12", "94-42=52")
    ,@("25-22=3", "6+33=39")
    ,@("13+31=44", "51-40=11")
    ,@("85+9=94", "11+21=32")
    ,@("37+29=66", "41+11=52")
    ,@("51+48=99", "36-0=36")
    ,@("87-27=60", "47+50=97")
    ,@("74-38=36", "47-5=42")
    ,@("40-14=26", "59+40=99")
    ,@("81+2=83", "31+67=98")
    ,@("28+38=66", "59-3=56")
    ,@("14+29=43", "17+18=35")
    ,@("41+23=64", "14+26=40")
    ,@("63-32=31", "12+5=17")
    ,@("60-36=24", "99-30=69")
    ,@("53+27=80", "6+29=35")
    ,@("34+55=89", "93-4=89")
    ,@("35+28=63", "1+62=63")
    ,@("80-34=46", "93-68=25")
    ,@("19+62=81", "52+19=71")
    ,@("35+27=62", "6+9=15")
    ,@("89-84=5", "67+29=96")
    ,@("90-17=73", "56+25=81")
    ,@("92-6=86", "83-52=31")
    ,@("60-9=51", "41+17=58")
    ,@("96-71=25", "69-4=65")
    ,@("42+37=79", "5+83=88")
    ,@("52+20=72", "43+16=59")
    ,@("32-7=25", "51-27=24")
    ,@("78-61=17", "34+52=86")
    ,@("38-15=23", "55-29=26")
    ,@("50+45=95", "82-73=9")
    ,@("75+4=79", "29-7=22")
    ,@("26+60=86", "53-52=1")
    ,@("83+10=93", "33+21=54")
    ,@("84-45=39", "71-3=68")
    ,@("83-38=45", "98-82=16")
    ,@("45-33=12", "8+89=97")
    ,@("97-19=78", "15+0=15")
    ,@("99-68=31", "82-19=63")
    ,@("65-2=63", "68+28=96")
    ,@("20+61=81", "38+42=80")
    ,@("96-84=12", "58-18=40")
    ,@("78-73=5", "9-3=6")
    ,@("87-68=19", "44+27=71")
    ,@("73-7=66", "3+48=51")
    ,@("40-33=7", "94-17=77")
    ,@("31+65=96", "52-50=2")
    ,@("24-15=9", "52+0=52")
    ,@("74+6=80", "54+35=89")
    ,@("20+62=82", "93-27=66")
    ,@("63-61=2", "72-53=19")
    ,@("74+12=86", "47+37=84")
    ,@("83+16=99", "18+34=52")
    ,@("13+33=46", "81-56=25")
    ,@("39+42=81", "43-5=38")
    ,@("18+55=73", "51+25=76")
    ,@("34+0=34", "88+11=99")
    ,@("66-24=42", "4+14=18")
    ,@("71-55=16", "49+22=71")
    ,@("0+9=9", "22+46=68")
    ,@("52-46=6", "77-49=28")
    ,@("48+3=51", "56+30=86")
    ,@("11+5=16", "80-45=35")
    ,@("22+30=52", "25+36=61")
    ,@("23+0=23", "18+41=59")
    ,@("21+66=87", "80-2=78")
    ,@("9+54=63", "53-2=51")
    ,@("12+8=20", "26+27=53")
    ,@("29-28=1", "17+78=95")
    ,@("68+18=86", "19+26=45")
    ,@("2+71=73", "88-37=51")
    ,@("9+72=81", "59-5=54")
    ,@("63+11=74", "45+12=57")
    ,@("84-51=33", "4+75=79")
    ,@("43+18=61", "49-25=24")
    ,@("92-65=27", "73-59=14")
    ,@("77-26=51", "34-18=16")
    ,@("11+52=63", "51-33=18")
    ,@("98-23=75", "35+22=57")
    ,@("33-14=19", "6+0=6")
    ,@("49-13=36", "92-40=52")
    ,@("2+13=15", "42+35=77")
    ,@("73-40=33", "44+47=91")
    ,@("74-34=40", "82-28=54")
    ,@("6+18=24", "6+45=51")
    ,@("34+54=88", "28-5=23")
    ,@("0+38=38", "59-55=4")
    ,@("3+9=12", "23+11=34")
    ,@("63-27=36", "23+51=74")
    ,@("66-42=24", "15-1=14")
    ,@("55+10=65", "20+7=27")
    ,@("12+59=71", "66-16=50")
    ,@("65+15=80", "67+18=85")
    ,@("88-28=60", "37+46=83")
    ,@("53+34=87", "70-57=13")
    ,@("31+40=71", "40-19=21")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done: replaced $($pairs.Count) text values"